# Update the OVS sheet: change the %-difference formula in column F so it
# uses an absolute-value comparison between the measured (B) and computed
# (E) readings instead of a plain ratio.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OVS")

$ws.Range("F2").Formula = "=ABS(E2-B2)/B2*100"
$ws.Range("F3").Formula = "=ABS(E3-B3)/B3*100"
$ws.Range("F4:F52").Formula = "=ABS(E4-B4)/B4*100"

# Make "OVS" the active sheet/tab and leave the selection on F3, matching
# where the user was working when the workbook was saved.
$ws.Activate()
$ws.Range("F3").Select()
